$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17: only column D (report cutoff date/time) changes
$ws.Range("D2:D17").Value = 45951.305277777778

# Rows 18-40: station (A), terminal (B), last-charge-end (C), and cutoff (D) all refreshed
$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "406号直流"
$ws.Range("C18").Value = 45949.623159722221
$ws.Range("D18").Value = 45951.305277777778

$ws.Range("A19").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B19").Value = "903号直流"
$ws.Range("C19").Value = 45950.249791666669
$ws.Range("D19").Value = 45951.305277777778

$ws.Range("A20").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B20").Value = "306号直流"
$ws.Range("C20").Value = 45950.272002314814
$ws.Range("D20").Value = 45951.305277777778

$ws.Range("A21").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B21").Value = "201号直流"
$ws.Range("C21").Value = 45950.436400462961
$ws.Range("D21").Value = 45951.305277777778

$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "603号直流"
$ws.Range("C22").Value = 45950.523645833331
$ws.Range("D22").Value = 45951.305277777778

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "405号直流"
$ws.Range("C23").Value = 45950.545925925922
$ws.Range("D23").Value = 45951.305277777778

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "505号直流"
$ws.Range("C24").Value = 45950.54991898148
$ws.Range("D24").Value = 45951.305277777778

$ws.Range("A25").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B25").Value = "109号直流"
$ws.Range("C25").Value = 45950.557083333333
$ws.Range("D25").Value = 45951.305277777778

$ws.Range("A26").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B26").Value = "002A号直流"
$ws.Range("C26").Value = 45950.560543981483
$ws.Range("D26").Value = 45951.305277777778

$ws.Range("A27").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value = "101号直流"
$ws.Range("C27").Value = 45950.574606481481
$ws.Range("D27").Value = 45951.305277777778

$ws.Range("A28").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B28").Value = "402号直流"
$ws.Range("C28").Value = 45950.574618055558
$ws.Range("D28").Value = 45951.305277777778

$ws.Range("A29").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B29").Value = "905号直流"
$ws.Range("C29").Value = 45950.593113425923
$ws.Range("D29").Value = 45951.305277777778

$ws.Range("A30").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B30").Value = "105号直流"
$ws.Range("C30").Value = 45950.595810185187
$ws.Range("D30").Value = 45951.305277777778

$ws.Range("A31").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B31").Value = "902号直流"
$ws.Range("C31").Value = 45950.603229166663
$ws.Range("D31").Value = 45951.305277777778

$ws.Range("A32").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B32").Value = "107号直流"
$ws.Range("C32").Value = 45950.605254629627
$ws.Range("D32").Value = 45951.305277777778

$ws.Range("A33").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B33").Value = "304号直流"
$ws.Range("C33").Value = 45950.628530092596
$ws.Range("D33").Value = 45951.305277777778

$ws.Range("A34").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B34").Value = "203号直流"
$ws.Range("C34").Value = 45950.655092592591
$ws.Range("D34").Value = 45951.305277777778

$ws.Range("A35").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B35").Value = "202号直流"
$ws.Range("C35").Value = 45950.672453703701
$ws.Range("D35").Value = 45951.305277777778

$ws.Range("A36").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B36").Value = "102号直流"
$ws.Range("C36").Value = 45950.700289351851
$ws.Range("D36").Value = 45951.305277777778

$ws.Range("A37").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B37").Value = "310号直流"
$ws.Range("C37").Value = 45950.730486111112
$ws.Range("D37").Value = 45951.305277777778

$ws.Range("A38").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B38").Value = "501号直流"
$ws.Range("C38").Value = 45950.736122685186
$ws.Range("D38").Value = 45951.305277777778

$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "103号直流"
$ws.Range("C39").Value = 45950.76898148148
$ws.Range("D39").Value = 45951.305277777778

$ws.Range("A40").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B40").Value = "002B号直流"
$ws.Range("C40").Value = 45950.77039351852
$ws.Range("D40").Value = 45951.305277777778

# Rows 41-42: source data exhausted -- clear to blank (styles retained)
$ws.Range("A41:D42").ClearContents()

# Selection moved from E12 to E10
$ws.Range("E10").Select()
